$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '30.170.43'
$ws.Cells.Item(2, 5).Value = '  -0.52%  '

$ws.Cells.Item(3, 4).Value = '1.914.23'
$ws.Cells.Item(3, 5).Value = '  -1.02%  '

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '0.9995'
$ws.Cells.Item(4, 5).Value = '  +0.01%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '0.7409'
$ws.Cells.Item(5, 5).Value = '  -2.41%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '243.94'
$ws.Cells.Item(6, 5).Value = '  -0.37%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.9988'
$ws.Cells.Item(7, 5).Value = '  +0.00%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.3126'
$ws.Cells.Item(8, 5).Value = '  -1.89%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '26.65'
$ws.Cells.Item(9, 5).Value = '  -3.61%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.06964'
$ws.Cells.Item(10, 5).Value = '  -0.51%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.7812'
$ws.Cells.Item(11, 5).Value = '  +0.32%  '

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.07982'
$ws.Cells.Item(12, 5).Value = '  -0.14%  '

$ws.Cells.Item(13, 4).Value = '1.907.63'
$ws.Cells.Item(13, 5).Value = '  -1.34%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '5.288'
$ws.Cells.Item(14, 5).Value = '  -1.22%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '92.36'
$ws.Cells.Item(15, 5).Value = '  -2.08%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '14.35'
$ws.Cells.Item(16, 5).Value = '  -0.53%  '

$ws.Cells.Item(17, 4).Value = '30.178.49'
$ws.Cells.Item(17, 5).Value = '  -0.50%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '5.903'
$ws.Cells.Item(18, 5).Value = '  +2.75%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '242.29'
$ws.Cells.Item(19, 5).Value = '  -4.21%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.000007827'
$ws.Cells.Item(20, 5).Value = '  -1.38%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '0.9992'
$ws.Cells.Item(21, 5).Value = '  +0.07%  '

$ws.Cells.Item(22, 4).Value = '2.134.00'
$ws.Cells.Item(22, 5).Value = '  -2.46%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.9994'
$ws.Cells.Item(23, 5).Value = '  +0.11%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '7.024'
$ws.Cells.Item(24, 5).Value = '  +5.09%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '9.407'
$ws.Cells.Item(25, 5).Value = '  -0.98%  '

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '168.07'
$ws.Cells.Item(26, 5).Value = '  +1.42%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '19.12'
$ws.Cells.Item(27, 5).Value = '  +0.76%  '

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '0.1284'
$ws.Cells.Item(28, 5).Value = '  -4.02%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.077'
$ws.Cells.Item(29, 5).Value = '  -5.04%  '

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '1.350'
$ws.Cells.Item(30, 5).Value = '  -0.90%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '1.547'
$ws.Cells.Item(31, 5).Value = '  +2.32%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '4.339'
$ws.Cells.Item(32, 5).Value = '  -1.02%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '4.106'
$ws.Cells.Item(33, 5).Value = '  -0.48%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.05170'
$ws.Cells.Item(34, 5).Value = '  +0.15%  '

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '1.306'
$ws.Cells.Item(35, 5).Value = '  +1.50%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.7474'
$ws.Cells.Item(36, 5).Value = '  -0.31%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.723'
$ws.Cells.Item(37, 5).Value = '  -1.62%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.01946'
$ws.Cells.Item(38, 5).Value = '  -0.47%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '2.798'
$ws.Cells.Item(39, 5).Value = '  -0.01%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '6.362'
$ws.Cells.Item(40, 5).Value = '  -0.85%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '75.09'
$ws.Cells.Item(41, 5).Value = '  -3.27%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.4506'
$ws.Cells.Item(42, 5).Value = '  +0.96%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '1.957'
$ws.Cells.Item(43, 5).Value = '  -0.43%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '7.886'
$ws.Cells.Item(44, 5).Value = '  +5.59%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '1.000'
$ws.Cells.Item(45, 5).Value = '  +0.10%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.8388'
$ws.Cells.Item(46, 5).Value = '  +0.70%  '

$ws.Cells.Item(47, 2).Value = 'Quant'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '101.56'
$ws.Cells.Item(47, 5).Value = '  +0.87%  '

$ws.Cells.Item(48, 2).Value = 'EnergySwap'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '9.917'
$ws.Cells.Item(48, 5).Value = '  +1.51%  '

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '37.43'
$ws.Cells.Item(49, 5).Value = '  +0.12%  '

$ws.Cells.Item(50, 4).Value = '2.037.57'
$ws.Cells.Item(50, 5).Value = '  -2.18%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '939.36'
$ws.Cells.Item(51, 5).Value = '  -4.58%  '
